$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column G: "matrixColorRelFormula" header + client WS call value ---

# Header cell G2 (style copied from F2's header style so the new column
# matches the look of the existing PredictiveDatasourceName header)
$ws.Range("F2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Value = "matrixColorRelFormula"

# G1 keeps the same (empty) banner-row styling as F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

# Data row value
$ws.Range("G3").Value = "LIB_ISP.updateMatrixColorRelation"

# Column G width, sized (bestFit-ish) for the header text
$ws.Columns.Item(7).ColumnWidth = 31.75

# Row 2 grows to fit the new, taller wrapped header text
$ws.Rows.Item(2).RowHeight = 45

# Keep the selection / active cell on the new last data cell, like the source
$ws.Range("G3").Select() | Out-Null

$excel.CutCopyMode = 0
